$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '29.935.10'
$ws.Range("E2").Value = '  +0.85%  '
$ws.Range("D3").Value = '1.632.77'
$ws.Range("E3").Value = '  +1.70%  '
$ws.Range("E4").Value = '  -0.14%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '214.53'
$ws.Range("E5").Value = '  +0.81%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '0.518'
$ws.Range("E6").Value = '  +0.17%  '
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.999'
$ws.Range("E7").Value = '  -0.14%  '
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '28.66'
$ws.Range("E8").Value = '  +2.22%  '
$ws.Range("E9").Value = '  +2.16%  '
$ws.Range("E10").Value = '  +0.88%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.0913'
$ws.Range("E11").Value = '  +0.28%  '
$ws.Range("D12").Value = '1.865.73'
$ws.Range("E12").Value = '  +1.62%  '
$ws.Range("D13").Value = '1.637.04'
$ws.Range("E13").Value = '  +1.88%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '0.564'
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '9.31'
$ws.Range("E15").Value = '  +17.95%  '
$ws.Range("B16").Value = 'WrappedBTC'
$ws.Range("C16").Value = 'https://coinranking.com/coin/x4WXHge-vvFY+wrappedbtc-wbtc'
$ws.Range("D16").Value = '29.925.84'
$ws.Range("E16").Value = '  +0.76%  '
$ws.Range("B17").Value = 'Polkadot'
$ws.Range("C17").Value = 'https://coinranking.com/coin/25W7FG7om+polkadot-dot'
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '3.86'
$ws.Range("E17").Value = '  +2.63%  '
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '64.09'
$ws.Range("E18").Value = '  +0.06%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '242.89'
$ws.Range("E19").Value = '  +0.83%  '
$ws.Range("D20").Value = '0.0₃0701'
$ws.Range("E20").Value = '  +0.50%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '0.999'
$ws.Range("E21").Value = '  -0.10%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '9.83'
$ws.Range("E22").Value = '  +4.54%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '4.13'
$ws.Range("E23").Value = '  +2.52%  '
$ws.Range("E24").Value = '  +1.22%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '157.61'
$ws.Range("E25").Value = '  +1.72%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '15.51'
$ws.Range("E26").Value = '  +0.25%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '0.110'
$ws.Range("E27").Value = '  +1.19%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '6.58'
$ws.Range("E28").Value = '  +1.99%  '
$ws.Range("E29").Value = '  -0.18%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '0.0485'
$ws.Range("E30").Value = '  +1.08%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '1.12'
$ws.Range("E31").Value = '  +4.53%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '3.38'
$ws.Range("E32").Value = '  +4.43%  '
$ws.Range("E33").Value = '  -0.66%  '
$ws.Range("D34").Value = '1.423.49'
$ws.Range("E34").Value = '  -0.20%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '1.64'
$ws.Range("E35").Value = '  +4.39%  '
$ws.Range("E36").Value = '  -0.05%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '2.80'
$ws.Range("E37").Value = '  -3.62%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '2.29'
$ws.Range("E38").Value = '  -0.38%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '0.0170'
$ws.Range("E39").Value = '  +0.38%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '75.25'
$ws.Range("E40").Value = '  +13.54%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '0.550'
$ws.Range("E41").Value = '  +0.33%  '
$ws.Range("E42").Value = '  +2.05%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '0.827'
$ws.Range("E43").Value = '  +1.40%  '
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '0.0490'
$ws.Range("E44").Value = '  -1.25%  '
$ws.Range("B45").Value = 'BitcoinSV'
$ws.Range("C45").Value = 'https://coinranking.com/coin/VcMY11NONHSA0+bitcoinsv-bsv'
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '53.25'
$ws.Range("E45").Value = '  -6.16%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '0.999'
$ws.Range("E46").Value = '  -0.13%  '
$ws.Range("B47").Value = 'WEMIXToken'
$ws.Range("C47").Value = 'https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix'
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '1.01'
$ws.Range("E47").Value = '  +3.28%  '
$ws.Range("B48").Value = 'RocketPoolETH'
$ws.Range("C48").Value = 'https://coinranking.com/coin/QJZRUGyNI+rocketpooleth-reth'
$ws.Range("D48").Value = '1.773.95'
$ws.Range("E48").Value = '  +1.80%  '
$ws.Range("B49").Value = 'FraxShare'
$ws.Range("C49").Value = 'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs'
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '5.34'
$ws.Range("E49").Value = '  -0.77%  '
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '89.18'
$ws.Range("E50").Value = '  +2.89%  '
$ws.Range("D51").Value = '0.0₆0111'
$ws.Range("E51").Value = '  +5.93%  '
